$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 397, shifting the existing rows 397:472 down to 398:473.
$ws.Range("A397").EntireRow.Insert()

# Populate the newly inserted row 397 with the new record.
$ws.Cells.Item(397, 1).Value = 4
$ws.Cells.Item(397, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(397, 3).Value = "Los Lagos"
$ws.Cells.Item(397, 4).Value = 45211
$ws.Cells.Item(397, 5).Value = 10
$ws.Cells.Item(397, 6).Value = 100112021
$ws.Cells.Item(397, 7).Value = "Ají"
$ws.Cells.Item(397, 8).Value = "Inferno"
$ws.Cells.Item(397, 9).Value = "Primera"
$ws.Cells.Item(397, 10).Value = 80
$ws.Cells.Item(397, 11).Value = 40000
$ws.Cells.Item(397, 12).Value = 40000
$ws.Cells.Item(397, 13).Value = 40000
$ws.Cells.Item(397, 14).Value = "$/caja 12 kilos"
$ws.Cells.Item(397, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(397, 16).Value = 4000
$ws.Cells.Item(397, 17).Value = 12
$ws.Cells.Item(397, 18).Value = "Hortaliza"
